$d = $word.ActiveDocument

# Step 1: merge the "Compile the Test project..." runs via self find/replace
$d.Content.Find.Execute("Compile the " + [char]0x201C + "Test" + [char]0x201D + " project. That is, the thing you got from the git repository. You might want to change the project settings and configuration to suit your needs. I have it set exclusively to x64 for my personal needs, because my Visual Studio 2010 is bugged (x86 doesn" + [char]0x2019 + "t work).", $false, $false, $false, $false, $false, $true, 1, $false, "Compile the " + [char]0x201C + "Test" + [char]0x201D + " project. That is, the thing you got from the git repository. You might want to change the project settings and configuration to suit your needs. I have it set exclusively to x64 for my personal needs, because my Visual Studio 2010 is bugged (x86 doesn" + [char]0x2019 + "t work).", 2) | Out-Null

$rng = $d.Content
$found = $rng.Find.Execute("bugged (x86 doesn" + [char]0x2019 + "t work).", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output "Found: $found text=$($rng.Text)"
$rng.Collapse(0)
$rng.InsertAfter(" If you got the code pre-built, this step is very skippable.")
